$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Lexican"
$ws.Range("C1").Value = "Designation"
$ws.Range("D1").Value = "Connector"
$ws.Range("E1").Value = "Connector1"
$ws.Range("F1").Value = "Predicate"
$ws.Range("G1").Value = "Concept"
$ws.Range("H1").Value = "Subject"
$ws.Range("I1").Value = "predicate_designation"
$ws.Range("J1").Value = "Definition"
$ws.Range("K1").Value = "Note"

$ws.Range("A2").Value = "testCaseID_01"
$ws.Range("B2").Value = "Jack_test"
$ws.Range("C2").Value = "Concrete_test"
$ws.Range("D2").Value = "11 - a"
$ws.Range("E2").Value = "34 - is the"
$ws.Range("F2").Value = "2 - pure setting"
$ws.Range("G2").Value = "9 - bridge"
$ws.Range("H2").Value = "9 - bridge"
$ws.Range("I2").Value = "Algebra"
$ws.Range("J2").Value = "Measurement"
$ws.Range("K2").Value = "Calculation"
